$p = $ppt.ActivePresentation

# Insert a new slide at position 6 (before the existing "Questions??" slide),
# using the "Title and Content" layout (CustomLayout index 2 on slideMaster1 -
# the same layout used by the "Coding" / "Good Design" slides already in the
# deck) so the new slide gets a Title placeholder + a Content placeholder.
$newSlide = $p.Slides.Add(6, 2)

# Title placeholder (shape 1)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Good Unit Tests"

# Content placeholder (shape 2): "Atomic" followed by a trailing empty paragraph
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Atomic"
[void]$body.InsertAfter([char]13)
